# CoverLetter.docx edit: add cover-letter body paragraph about Zlien + tidy up
# paragraph-mark rFonts hints that Word drops on the next real save.
$d = $word.ActiveDocument

# 1) Divider line paragraph ("ـــ...") - the paragraph-mark rPr no longer carries the
#    rFonts hint="cs" (the run itself keeps it).
$divider = $d.Paragraphs(6)
$dividerXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p w:rsidR="00B84DA9" w:rsidRDefault="00B84DA9" w:rsidP="00B84DA9"><w:pPr><w:ind w:left="0"/><w:jc w:val="right"/><w:rPr><w:rtl/><w:lang w:bidi="ar-EG"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="ar-EG"/></w:rPr><w:t>ــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــــ</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$divider.Range.InsertXML($dividerXml)

# 2) The blank RTL paragraph right after the divider loses its rFonts hint="cs" too.
$blankRtl = $d.Paragraphs(7)
$blankRtlXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p w:rsidR="00B84DA9" w:rsidRDefault="00B84DA9" w:rsidP="00B84DA9"><w:pPr><w:ind w:left="0"/><w:rPr><w:rtl/><w:lang w:bidi="ar-EG"/></w:rPr></w:pPr></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$blankRtl.Range.InsertXML($blankRtlXml)

# 3) Closing paragraph ("For the last 2 years ... (EC3). ") gains a trailing space,
#    drops the stray bookmark run at its end, and is followed by the new "According to
#    my research ..." body paragraph, a blank paragraph, then a final paragraph that
#    now hosts the _GoBack bookmark.
$closing = $d.Paragraphs.Last
$closingXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p w:rsidR="005C641D" w:rsidRPr="00CE0B01" w:rsidRDefault="005130B8" w:rsidP="00B84DA9"><w:pPr><w:ind w:left="0"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="ar-EG"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="ar-EG"/></w:rPr><w:t xml:space="preserve">For the last 2 years, I have been working as </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="ar-EG"/></w:rPr><w:t>a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="ar-EG"/></w:rPr><w:t xml:space="preserve"> infrastructure and security engineer at Egyptian Cloud Computer Center (EC3).  </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="0"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="ar-EG"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>According to my research</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Zlien</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> company has </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>managed to increase its project load each quarter and strengthen its already laudable position in the industry</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/><w:color w:val="000000"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> of Marketing, Advertising and Computer Software.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="0"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="ar-EG"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:ind w:left="0"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="ar-EG"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$closing.Range.InsertXML($closingXml)
